$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("届出")

# Add a light-cyan background fill to the attendance-table header (rows 6-7)
$ws.Range("A6:L7").Interior.Color = 16777164

# Merge the header cells of the table so the header reads cleanly
$ws.Range("A6:B6").Merge()
$ws.Range("D6:F6").Merge()
$ws.Range("G6:I6").Merge()
$ws.Range("C6:C7").Merge()
$ws.Range("J6:J7").Merge()
$ws.Range("K6:K7").Merge()
$ws.Range("L6:L7").Merge()

# Normalize the body/label font from the old Japanese Gothic font to Calibri
$ws.Range("A3:D4").Font.Name = "Calibri"
$ws.Range("A8:L47").Font.Name = "Calibri"
